# Rename worksheets (in sheet order) to reflect the rerun/summarise of
# models without urban landuse.
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ48304744",
    "summ50807931",
    "summ53404967",
    "summ56012905",
    "summ58687873",
    "summ02108541",
    "summ05082165",
    "summ07861240",
    "summ10682892"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}
